$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 352, shifting all following
# rows (352-418) down to (354-420).
$ws.Range("352:353").EntireRow.Insert()

# New row 352
$ws.Range("A352").Value = 8
$ws.Range("B352").Value = "Terminal La Palmera de La Serena"
$ws.Range("C352").Value = "Coquimbo"
$ws.Range("D352").Value = 44476
$ws.Range("E352").Value = 4
$ws.Range("F352").Value = 100112004
$ws.Range("G352").Value = "Cebolla"
$ws.Range("H352").Value = "Sin especificar"
$ws.Range("I352").Value = "1a (guarda)"
$ws.Range("J352").Value = 2560
$ws.Range("K352").Value = 4800
$ws.Range("L352").Value = 5000
$ws.Range("M352").Value = 4900
$ws.Range("N352").Value = "`$/malla 16 kilos"
$ws.Range("O352").Value = "Región de O'Higgins"
$ws.Range("P352").Value = 306
$ws.Range("Q352").Value = 16
$ws.Range("R352").Value = "Hortaliza"

# New row 353
$ws.Range("A353").Value = 8
$ws.Range("B353").Value = "Terminal La Palmera de La Serena"
$ws.Range("C353").Value = "Coquimbo"
$ws.Range("D353").Value = 44476
$ws.Range("E353").Value = 4
$ws.Range("F353").Value = 100112004
$ws.Range("G353").Value = "Cebolla"
$ws.Range("H353").Value = "Sin especificar"
$ws.Range("I353").Value = "2a (guarda)"
$ws.Range("J353").Value = 1460
$ws.Range("K353").Value = 4500
$ws.Range("L353").Value = 4600
$ws.Range("M353").Value = 4550
$ws.Range("N353").Value = "`$/malla 16 kilos"
$ws.Range("O353").Value = "Región de O'Higgins"
$ws.Range("P353").Value = 284
$ws.Range("Q353").Value = 16
$ws.Range("R353").Value = "Hortaliza"
